$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "60.845.63"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -2.81%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.353.30"
$ws.Range("D3").Style = "Normal"
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "566.03"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.11%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "146.65"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.32%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.480"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.30%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "7.91"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.02%  "
$ws.Range("E10").Value = "  -1.17%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.415"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.99%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "3.929.55"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.37%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "27.89"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.92%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.361.14"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.64%  "
$ws.Range("E16").Value = "  -1.24%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "60.908.52"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -2.80%  "
$ws.Range("E18").Value = "  -1.22%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "14.40"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.99%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "8.88"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.78%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "375.78"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.93%  "
$ws.Range("E22").Value = "  -0.10%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "74.87"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.39%  "
$ws.Range("E24").Value = "  -0.01%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.502.54"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.29%  "
$ws.Range("E26").Value = "  -6.18%  "
$ws.Range("E27").Value = "  -4.15%  "
$ws.Range("E28").Value = "  -0.08%  "
$ws.Range("E29").Value = "  -2.83%  "
$ws.Range("B30").Value = "USDe"
$ws.Range("C30").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.00"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.01%  "
$ws.Range("B31").Value = "PancakeSwap"
$ws.Range("C31").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.08"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.95%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.68"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -3.55%  "
$ws.Range("E33").Value = "  -2.25%  "
$ws.Range("E34").Value = "  -1.41%  "
$ws.Range("E35").Value = "  +0.46%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "168.81"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.92%  "
$ws.Range("E37").Value = "  -4.11%  "
$ws.Range("E38").Value = "  -2.26%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "28.91"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -9.74%  "
$ws.Range("E40").Value = "  -2.30%  "
$ws.Range("E41").Value = "  -2.70%  "
$ws.Range("E42").Value = "  -3.41%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "4.29"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.01%  "
$ws.Range("E44").Value = "  -3.22%  "
$ws.Range("E45").Value = "  -5.22%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.488.25"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.72%  "
$ws.Range("B47").Value = "Cosmos"
$ws.Range("C47").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "6.66"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -3.36%  "
$ws.Range("B48").Value = "InjectiveProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "22.55"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.63%  "
$ws.Range("E49").Value = "  -0.02%  "
$ws.Range("E50").Value = "  -2.22%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.812"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.20%  "
